# "Natmi following Dr Hou advice" — update the NATMI ligand-receptor
# metrics for rows 2-26 (columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T) to the
# recomputed values. Columns A-D, F and L are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1-based column numbers for E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$colIndices = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

# New values for rows 2..26, in the same column order as $colIndices.
# Stored as strings and cast to [double] so values using scientific
# notation (e.g. "...E-05") parse correctly.
$newValues = @(
    @("3","22.481209","67.44362700000001","0.1656226259370683","0.166106832923046","3","0.1757273333333333","0.527182","0.00403052093865001","0.004035374510811824","3.950562907679334","35.55506616911401","0.0006675454617535522","0.000670303279649338"),
    @("3","22.481209","67.44362700000001","0.1656226259370683","0.166106832923046","3","0.362774","1.088322","0.008320664607276909","0.008330684390506021","8.155598113766001","73.40038302389401","0.001378090321798827","0.001383783600188411"),
    @("3","22.481209","67.44362700000001","0.1656226259370683","0.166106832923046","3","18.888837","56.666511","0.4332385383145499","0.4337602461882952","424.643892363933","3821.795031275397","0.07175410437279295","0.07205054074225845"),
    @("3","22.481209","67.44362700000001","0.1656226259370683","0.166106832923046","3","24.01450566666666","72.04351699999999","0.5508020072052686","0.5514652855580013","539.8751209240177","4858.876088316159","0.09122527480474463","0.09160215205104277"),
    @("3","22.481209","67.44362700000001","0.1656226259370683","0.166106832923046","2","0.1573175","0.314635","0.003608268934254619","0.002408409352385472","3.5366875968575","21.220125581145","0.000597610975978397","0.000400053249906995"),
    @("3","29.50180766666667","88.50542300000001","0.2173444878184117","0.2179799065528387","3","0.1757273333333333","0.527182","0.00403052093865001","0.004035374510811824","5.184273989776223","46.65846590798601","0.0008760115090522705","0.0008796305587724686"),
    @("3","29.50180766666667","88.50542300000001","0.2173444878184117","0.2179799065528387","3","0.362774","1.088322","0.008320664607276909","0.008330684390506021","10.70248877446733","96.32239897020601","0.001808450587377386","0.001815921804963695"),
    @("3","29.50180766666667","88.50542300000001","0.2173444878184117","0.2179799065528387","3","18.888837","56.666511","0.4332385383145499","0.4337602461882952","557.2548362210169","5015.293525989154","0.09416200821317319","0.0945510179304609"),
    @("3","29.50180766666667","88.50542300000001","0.2173444878184117","0.2179799065528387","3","24.01450566666666","72.04351699999999","0.5508020072052686","0.5514652855580013","708.4713273880767","6376.241946492691","0.1197137801453822","0.1202083514130677"),
    @("3","29.50180766666667","88.50542300000001","0.2173444878184117","0.2179799065528387","2","0.1573175","0.314635","0.003608268934254619","0.002408409352385472","4.641150627600833","27.846903765605","0.0007842373634266564","0.0005249848455739681"),
    @("3","42.765269","128.295807","0.3150585073376215","0.3159795984589671","3","0.1757273333333333","0.527182","0.00403052093865001","0.004035374510811824","7.515026680652666","67.63524012587401","0.001269849910724101","0.001275096017557871"),
    @("3","42.765269","128.295807","0.3150585073376215","0.3159795984589671","3","0.362774","1.088322","0.008320664607276909","0.008330684390506021","15.514127696206","139.627149265854","0.00262149617122564","0.002632326308600477"),
    @("3","42.765269","128.295807","0.3150585073376215","0.3159795984589671","3","18.888837","56.666511","0.4332385383145499","0.4337602461882952","807.7861954021529","7270.075758619377","0.136495487202515","0.1370593884180402"),
    @("3","42.765269","128.295807","0.3150585073376215","0.3159795984589671","3","24.01450566666666","72.04351699999999","0.5508020072052686","0.5514652855580013","1026.986794737024","9242.881152633217","0.1735348582286577","0.1742517794946769"),
    @("3","42.765269","128.295807","0.3150585073376215","0.3159795984589671","2","0.1573175","0.314635","0.003608268934254619","0.002408409352385472","6.727725205907499","40.366351235445","0.00113681582449897","0.0007610082200915824"),
    @("3","39.80222300000001","119.406669","0.2932292783449848","0.2940865504976542","3","0.1757273333333333","0.527182","0.00403052093865001","0.004035374510811824","6.994338508528668","62.94904657675801","0.001181866746194693","0.001186749369850808"),
    @("3","39.80222300000001","119.406669","0.2932292783449848","0.2940865504976542","3","0.362774","1.088322","0.008320664607276909","0.008330684390506021","14.439211646602","129.952904819418","0.002439862478142465","0.002449942235688569"),
    @("3","39.80222300000001","119.406669","0.2932292783449848","0.2940865504976542","3","18.888837","56.666511","0.4332385383145499","0.4337602461882952","751.817702484651","6766.359322361859","0.1270382239412115","0.127563054544529"),
    @("3","39.80222300000001","119.406669","0.2932292783449848","0.2940865504976542","3","24.01450566666666","72.04351699999999","0.5508020072052686","0.5514652855580013","955.8307097794303","8602.476388014873","0.16151127508377","0.1621785235489565"),
    @("3","39.80222300000001","119.406669","0.2932292783449848","0.2940865504976542","2","0.1573175","0.314635","0.003608268934254619","0.002408409352385472","6.261586216802501","37.569517300815","0.001058050095666109","0.0007082807986293328"),
    @("2","1.1870385","2.374077","0.008745100561913671","0.005847111567493934","3","0.1757273333333333","0.527182","0.00403052093865001","0.004035374510811824","0.208595110169","1.251570661014","3.524731092539302E-05","2.359528498133799E-05"),
    @("2","1.1870385","2.374077","0.008745100561913671","0.005847111567493934","3","0.362774","1.088322","0.008320664607276909","0.008330684390506021","0.430626704799","2.583760228794","7.276504873259249E-05","4.871044106486891E-05"),
    @("2","1.1870385","2.374077","0.008745100561913671","0.005847111567493934","3","18.888837","56.666511","0.4332385383145499","0.4337602461882952","22.4217767392245","134.530660435347","0.003788714584857228","0.002536244553006598"),
    @("2","1.1870385","2.374077","0.008745100561913671","0.005847111567493934","3","24.01450566666666","72.04351699999999","0.5508020072052686","0.5514652855580013","28.5061427848015","171.036856708809","0.004816818942713972","0.003224479050257535"),
    @("2","1.1870385","2.374077","0.008745100561913671","0.005847111567493934","2","0.1573175","0.314635","0.003608268934254619","0.002408409352385472","0.18674192922375","0.7469677168950001","3.155467468448571E-05","1.408223818359367E-05")
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = 2 + $i
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $colIndices.Count; $j++) {
        $col = $colIndices[$j]
        $ws.Cells.Item($row, $col).Value = [double]$rowVals[$j]
    }
}

Write-Output "Updated rows 2-26 for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T"
